# Add a new "investigators" column to the "borehole" sheet, between the
# existing "notes" (P) and "funding" (Q) columns, and update the "funding"
# comment wording (gastaldello2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("borehole")

# 1. Insert a new column at Q (17) -- this shifts the old Q column
#    ("funding", values + style + column-width bucket) one to the right,
#    to R (18).
$ws.Columns.Item(17).Insert()

# 2. Header value/label for the new column.
$ws.Range("Q1").Value = "investigators"

# 3. Column widths: P (16) keeps its width; the new investigators column
#    (17) gets a wider custom width; funding's old width carries over to
#    R (18) automatically from the Insert() above.
$ws.Range("Q1").ColumnWidth = 15.5

# 4. Comments.
#    Column-insert does not relocate the legacy VML comment anchors, so the
#    comment that used to describe "funding" is still sitting on the literal
#    Q1 cell -- which is exactly where the new "investigators" header is now.
#    Rewrite its text in place (this keeps the original comment author).
$investigatorsComment = "[string] investigators`nNames of people and/or agencies who performed the work, as a pipe-delimited list. Each entry should be in the format {person} ({agencies}) [{notes}], where either person or at least one (semicolon-delimited) agencies is required.`nconstraints:`n  - pattern: [^\s]+( [^\s]+)*"
[void]$ws.Range("Q1").Comment.Text($investigatorsComment)

#    Add a fresh comment for the relocated "funding" header at R1, with the
#    updated wording.
$fundingComment = "[string] funding`nFunding sources as a pipe-delimited list. Each entry should be in the format {funder} [{rorid}] > {award} [{number}] ({url}), where only the funder is required and rorid is the funder's ROR (https://ror.org) ID (e.g. 01jtrvx49).`nconstraints:`n  - pattern: [^\s]+( [^\s]+)*"
$ws.Range("R1").AddComment($fundingComment)

# 5. Conditional-formatting formulas reference the last column of the row
#    (previously Q == 17th column) to detect "all blank" rows; bump them to
#    R == 18th column now that there's one more column.
$ranges = @("A2:A1048576", "B2:B1048576", "D2:D1048576", "E2:E1048576", "F2:F1048576")
foreach ($addr in $ranges) {
    $rng = $ws.Range($addr)
    $fcs = $rng.FormatConditions()
    $fc = $fcs.Item(1)
    $formula = $fc.Formula1()
    $formula = $formula.Replace('$A2:$Q2', '$A2:$R2')
    $formula = $formula.Replace('<> 17', '<> 18')
    $fc.Formula1 = $formula
}
